$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("trial_total") values for rows 2-42 decrease by 106
# (228->122, 229->123, ... 268->162), so that n distractor = n targets.
for ($row = 2; $row -le 42; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $old = $cell.Value2
    $cell.Value2 = $old - 106
}
